$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (price & 1h volume change)
# Swap of row 43/44 (Cronos <-> TrustWalletToken) is included as plain cell updates.

$ws.Range("D2").Value = "37.448.56"
$ws.Range("E2").Value = "  -0.67%  "
$ws.Range("D3").Value = "2.067.28"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.630"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.32%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.35"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.97%  "
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0778"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("D13").Value = "2.375.29"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("E14").Value = "  -0.33%  "
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").Value = "2.068.54"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "37.423.49"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.44"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.94"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.97%  "
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  +0.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.56%  "
$ws.Range("E26").Value = "  +6.78%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.57"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.133"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.37"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.01%  "
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.60"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.48"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.32%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("E39").Value = "  -1.30%  "
$ws.Range("E40").Value = "  +6.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0953"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.02%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.20"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.91%  "
$ws.Range("D45").Value = "1.478.90"
$ws.Range("E45").Value = "  +2.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "16.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.21%  "
$ws.Range("E47").Value = "  -1.60%  "
$ws.Range("E48").Value = "  -1.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.05%  "
$ws.Range("D51").Value = "2.258.92"
$ws.Range("E51").Value = "  -0.29%  "
